$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3464964993005633
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.1529057820181812
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 2.652525447291612

$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 5.488907176552729

$ws.Range("B4").Value = 0.7287194209349384
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 16.98373111632243
$ws.Range("E4").Value = 6.48142807727062
$ws.Range("G4").Value = 25.8471150734268

$ws.Range("B5").Value = 3.182878228561681
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 5.488907176552729

$ws.Range("B6").Value = 3.182878228561681
$ws.Range("C6").Value = 1.65323645889881
$ws.Range("D6").Value = 0.7127328510149897
$ws.Range("E6").Value = 0.4998867070740569
$ws.Range("G6").Value = 6.048734245549538
